$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix an error in Admittance matrix creation: values in column C (rows 2 and 5)
# were off by a factor of 10 and need correcting.
$ws.Range("C2").Value = 0.054030000000000002
$ws.Range("C5").Value = 0.056950000000000001

# Update the active selection to reflect where the editing left off.
$ws.Range("D17").Select()
